$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures for this run.
# Cells are stored as text (like "327.48" or "3.25%"), so force a text
# number format before assigning the value to avoid Excel auto-converting
# them into numeric/percentage values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.25%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.61%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.886"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "13.20%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08021"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.20%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.586"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.36%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.719"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.929"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.22%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.942"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.45%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9423"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1257"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.19%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1967"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.58%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.861"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "35.08%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09134"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.37%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03566"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "6.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09647"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.15%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001300"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-7.04%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006124"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.86%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.366"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.79%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.29%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1433"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.94%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2412"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.52%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04409"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.99%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.57%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004366"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.49%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001143"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-13.87%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02416"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.54%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05271"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.33%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007476"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.81%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1419"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.46%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008710"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.09%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01068"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "29.74%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006854"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.86%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.56%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003152"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.21%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001425"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-15.62%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.56%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.56%"
